$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: add P1=14, Q1=15
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(2, 2).Value = 3.46677384630982
$ws.Cells.Item(2, 3).Value = 1.019169272890139
$ws.Cells.Item(2, 4).Value = 0.04824190763466873
$ws.Cells.Item(2, 5).Value = 1.35169962688731
$ws.Cells.Item(2, 6).Value = 0.5084535737087279
$ws.Cells.Item(2, 7).Value = 0.0007896385517037115
$ws.Cells.Item(2, 8).Value = 0.01092094002687105
$ws.Cells.Item(2, 9).Value = 0.003536527785124033
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 1.435151477765999

$ws.Cells.Item(3, 2).Value = 3.02011174815118
$ws.Cells.Item(3, 3).Value = 0.8990605546520669
$ws.Cells.Item(3, 4).Value = 0.0433835340809452
$ws.Cells.Item(3, 5).Value = 1.176740016824183
$ws.Cells.Item(3, 6).Value = 0.468187038062311
$ws.Cells.Item(3, 7).Value = 0.000793366441376707
$ws.Cells.Item(3, 8).Value = 0.00766975902495437
$ws.Cells.Item(3, 9).Value = 0.001965288275290966
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 1.348594257482659

$ws.Cells.Item(4, 2).Value = 2.745669458445775
$ws.Cells.Item(4, 3).Value = 0.8258832939075091
$ws.Cells.Item(4, 4).Value = 0.04041041276209967
$ws.Cells.Item(4, 5).Value = 1.069666178277558
$ws.Cells.Item(4, 6).Value = 0.4441198563540496
$ws.Cells.Item(4, 7).Value = 0.0007957272508090085
$ws.Cells.Item(4, 8).Value = 0.005904022270549425
$ws.Cells.Item(4, 9).Value = 0.001270814676121557
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 1.297464469391542

$ws.Cells.Item(5, 2).Value = 2.633730900298815
$ws.Cells.Item(5, 3).Value = 0.7977722387358028
$ws.Cells.Item(5, 4).Value = 0.03928077046956702
$ws.Cells.Item(5, 5).Value = 1.02608314129057
$ws.Cells.Item(5, 6).Value = 0.4336700664951891
$ws.Cells.Item(5, 7).Value = 0.0007967124296093884
$ws.Cells.Item(5, 8).Value = 0.00523637476029204
$ws.Cells.Item(5, 9).Value = 0.001110499418202515
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 1.274436995028537

$ws.Cells.Item(6, 2).Value = 2.615103766652396
$ws.Cells.Item(6, 3).Value = 0.7950489896263662
$ws.Cells.Item(6, 4).Value = 0.03919122151063448
$ws.Cells.Item(6, 5).Value = 1.018830727099584
$ws.Cells.Item(6, 6).Value = 0.4309748955527013
$ws.Cells.Item(6, 7).Value = 0.0007968828780179578
$ws.Cells.Item(6, 8).Value = 0.005126972417599962
$ws.Cells.Item(6, 9).Value = 0.001168165231651308
$ws.Cells.Item(6, 15).Value = 0
$ws.Cells.Item(6, 16).Value = 0
$ws.Cells.Item(6, 17).Value = 1.267397047251904

$ws.Cells.Item(7, 2).Value = 2.74406360535454
$ws.Cells.Item(7, 3).Value = 0.8308092758339853
$ws.Cells.Item(7, 4).Value = 0.04066341026042863
$ws.Cells.Item(7, 5).Value = 1.069026406947373
$ws.Cells.Item(7, 6).Value = 0.4413300423773592
$ws.Cells.Item(7, 7).Value = 0.0007957557770997617
$ws.Cells.Item(7, 8).Value = 0.005890106154369357
$ws.Cells.Item(7, 9).Value = 0.001472853056091239
$ws.Cells.Item(7, 15).Value = 0
$ws.Cells.Item(7, 16).Value = 0
$ws.Cells.Item(7, 17).Value = 1.28828808501018

$ws.Cells.Item(8, 2).Value = 3.312653545137607
$ws.Cells.Item(8, 3).Value = 0.9847298999007421
$ws.Cells.Item(8, 4).Value = 0.04692401796475565
$ws.Cells.Item(8, 5).Value = 1.291216389130142
$ws.Cells.Item(8, 6).Value = 0.4909136236972671
$ws.Cells.Item(8, 7).Value = 0.0007909279253908628
$ws.Cells.Item(8, 8).Value = 0.009741767389663553
$ws.Cells.Item(8, 9).Value = 0.003174773938467546
$ws.Cells.Item(8, 15).Value = 0
$ws.Cells.Item(8, 16).Value = 0
$ws.Cells.Item(8, 17).Value = 1.393125671300481

$ws.Cells.Item(9, 2).Value = 4.427610317638539
$ws.Cells.Item(9, 3).Value = 1.282110248089054
$ws.Cells.Item(9, 4).Value = 0.0588090931990024
$ws.Cells.Item(9, 5).Value = 1.730915313986756
$ws.Cells.Item(9, 6).Value = 0.5988199823157387
$ws.Cells.Item(9, 7).Value = 0.0007819898167375172
$ws.Cells.Item(9, 8).Value = 0.01935802986678037
$ws.Cells.Item(9, 9).Value = 0.008607690668569923
$ws.Cells.Item(9, 15).Value = 0
$ws.Cells.Item(9, 16).Value = 0
$ws.Cells.Item(9, 17).Value = 1.632718893241275

$ws.Cells.Item(10, 2).Value = 5.24749449069725
$ws.Cells.Item(10, 3).Value = 1.50644471872306
$ws.Cells.Item(10, 4).Value = 0.06925275944912102
$ws.Cells.Item(10, 5).Value = 1.955498156088424
$ws.Cells.Item(10, 6).Value = 0.6683724546327312
$ws.Cells.Item(10, 7).Value = 0.0007759177289413412
$ws.Cells.Item(10, 8).Value = 0.02725836382229918
$ws.Cells.Item(10, 9).Value = 0.0143664161419883
$ws.Cells.Item(10, 15).Value = 0
$ws.Cells.Item(10, 16).Value = 0
$ws.Cells.Item(10, 17).Value = 1.775547518866148

$ws.Cells.Item(11, 2).Value = 5.612142301903873
$ws.Cells.Item(11, 3).Value = 1.610434802784823
$ws.Cells.Item(11, 4).Value = 0.08861588644301577
$ws.Cells.Item(11, 5).Value = 1.267337406619561
$ws.Cells.Item(11, 6).Value = 0.5871507068896307
$ws.Cells.Item(11, 7).Value = 0.0007746196692268896
$ws.Cells.Item(11, 8).Value = 0.04270868006982909
$ws.Cells.Item(11, 9).Value = 0.01643561128637661
$ws.Cells.Item(11, 15).Value = 0
$ws.Cells.Item(11, 16).Value = 0
$ws.Cells.Item(11, 17).Value = 1.474547141099578

$ws.Cells.Item(12, 2).Value = 5.74725555498668
$ws.Cells.Item(12, 3).Value = 1.642389293035194
$ws.Cells.Item(12, 4).Value = 0.1037805964393641
$ws.Cells.Item(12, 5).Value = 0.7713189373199754
$ws.Cells.Item(12, 6).Value = 0.5121260620343264
$ws.Cells.Item(12, 7).Value = 0.0007745905169331831
$ws.Cells.Item(12, 8).Value = 0.07854308646568597
$ws.Cells.Item(12, 9).Value = 0.01664005842228544
$ws.Cells.Item(12, 15).Value = 0
$ws.Cells.Item(12, 16).Value = 0
$ws.Cells.Item(12, 17).Value = 1.224683832186997

$ws.Cells.Item(13, 2).Value = 5.71206086016997
$ws.Cells.Item(13, 3).Value = 1.628934670659646
$ws.Cells.Item(13, 4).Value = 0.1168170525753141
$ws.Cells.Item(13, 5).Value = 0.3948377359290731
$ws.Cells.Item(13, 6).Value = 0.4334589186537912
$ws.Cells.Item(13, 7).Value = 0.0007755528712741501
$ws.Cells.Item(13, 8).Value = 0.1314041879758321
$ws.Cells.Item(13, 9).Value = 0.01567717074598907
$ws.Cells.Item(13, 15).Value = 0
$ws.Cells.Item(13, 16).Value = 0
$ws.Cells.Item(13, 17).Value = 0.9858811040282944

$ws.Cells.Item(14, 2).Value = 5.610402796626545
$ws.Cells.Item(14, 3).Value = 1.600782500645209
$ws.Cells.Item(14, 4).Value = 0.1251548963565625
$ws.Cells.Item(14, 5).Value = 0.2040035724406977
$ws.Cells.Item(14, 6).Value = 0.3772055239524477
$ws.Cells.Item(14, 7).Value = 0.0007766576504858142
$ws.Cells.Item(14, 8).Value = 0.1788295655706946
$ws.Cells.Item(14, 9).Value = 0.01459786682350295
$ws.Cells.Item(14, 15).Value = 0
$ws.Cells.Item(14, 16).Value = 0
$ws.Cells.Item(14, 17).Value = 0.8260307867355579

$ws.Cells.Item(15, 2).Value = 5.549222564864976
$ws.Cells.Item(15, 3).Value = 1.586479300100507
$ws.Cells.Item(15, 4).Value = 0.1266680095887693
$ws.Cells.Item(15, 5).Value = 0.1664673984577796
$ws.Cells.Item(15, 6).Value = 0.361549071281118
$ws.Cells.Item(15, 7).Value = 0.0007771664842525621
$ws.Cells.Item(15, 8).Value = 0.1906857109769646
$ws.Cells.Item(15, 9).Value = 0.01418120435697556
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 0.7848470606577962

$ws.Cells.Item(16, 2).Value = 5.202401853822494
$ws.Cells.Item(16, 3).Value = 1.495615476278033
$ws.Cells.Item(16, 4).Value = 0.1191748546277722
$ws.Cells.Item(16, 5).Value = 0.1619581450697751
$ws.Cells.Item(16, 6).Value = 0.3473887037929586
$ws.Cells.Item(16, 7).Value = 0.0007795445118418047
$ws.Cells.Item(16, 8).Value = 0.175491235802042
$ws.Cells.Item(16, 9).Value = 0.01192358401845794
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 0.7706965671103063

$ws.Cells.Item(17, 2).Value = 4.991035451832488
$ws.Cells.Item(17, 3).Value = 1.441267988981963
$ws.Cells.Item(17, 4).Value = 0.1089952125854552
$ws.Cells.Item(17, 5).Value = 0.2459753982704243
$ws.Cells.Item(17, 6).Value = 0.365796307107324
$ws.Cells.Item(17, 7).Value = 0.000780802229109395
$ws.Cells.Item(17, 8).Value = 0.1367689001963726
$ws.Cells.Item(17, 9).Value = 0.01081278822335374
$ws.Cells.Item(17, 15).Value = 0
$ws.Cells.Item(17, 16).Value = 0
$ws.Cells.Item(17, 17).Value = 0.8422055574621652

$ws.Cells.Item(18, 2).Value = 4.871728470037112
$ws.Cells.Item(18, 3).Value = 1.407347280808153
$ws.Cells.Item(18, 4).Value = 0.09563455489917061
$ws.Cells.Item(18, 5).Value = 0.4790218023010837
$ws.Cells.Item(18, 6).Value = 0.4181691898324047
$ws.Cells.Item(18, 7).Value = 0.0007811315644281613
$ws.Cells.Item(18, 8).Value = 0.08461431799794639
$ws.Cells.Item(18, 9).Value = 0.01026415122278657
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = 1.011170216096644

$ws.Cells.Item(19, 2).Value = 4.83500880465499
$ws.Cells.Item(19, 3).Value = 1.40339474351282
$ws.Cells.Item(19, 4).Value = 0.08223365033970964
$ws.Cells.Item(19, 5).Value = 0.9131096264680565
$ws.Cells.Item(19, 6).Value = 0.4940337535870114
$ws.Cells.Item(19, 7).Value = 0.0007805737522911431
$ws.Cells.Item(19, 8).Value = 0.0420758044946794
$ws.Cells.Item(19, 9).Value = 0.01072812603060136
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = 1.251779364725252

$ws.Cells.Item(20, 2).Value = 5.031627643970751
$ws.Cells.Item(20, 3).Value = 1.464437765624723
$ws.Cells.Item(20, 4).Value = 0.06744620694644254
$ws.Cells.Item(20, 5).Value = 1.891890550245648
$ws.Cells.Item(20, 6).Value = 0.6409981226065185
$ws.Cells.Item(20, 7).Value = 0.0007775358304290444
$ws.Cells.Item(20, 8).Value = 0.02500230937025849
$ws.Cells.Item(20, 9).Value = 0.01333659226125761
$ws.Cells.Item(20, 15).Value = 0
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = 1.707873429777266

$ws.Cells.Item(21, 2).Value = 5.663603663768583
$ws.Cells.Item(21, 3).Value = 1.636139845533137
$ws.Cells.Item(21, 4).Value = 0.07300782904698622
$ws.Cells.Item(21, 5).Value = 2.224275313971731
$ws.Cells.Item(21, 6).Value = 0.7181090159483858
$ws.Cells.Item(21, 7).Value = 0.0007726865873697398
$ws.Cells.Item(21, 8).Value = 0.03281830044685519
$ws.Cells.Item(21, 9).Value = 0.01854060229383236
$ws.Cells.Item(21, 15).Value = 0
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = 1.895197017569245

$ws.Cells.Item(22, 2).Value = 6.077296100170315
$ws.Cells.Item(22, 3).Value = 1.741715857843474
$ws.Cells.Item(22, 4).Value = 0.07705961331203781
$ws.Cells.Item(22, 5).Value = 2.390794747381719
$ws.Cells.Item(22, 6).Value = 0.7665380570823004
$ws.Cells.Item(22, 7).Value = 0.000769642109448282
$ws.Cells.Item(22, 8).Value = 0.03795447932891705
$ws.Cells.Item(22, 9).Value = 0.02207615676014552
$ws.Cells.Item(22, 15).Value = 0
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = 2.010982655735347

$ws.Cells.Item(23, 2).Value = 5.856533976777541
$ws.Cells.Item(23, 3).Value = 1.678693843975964
$ws.Cells.Item(23, 4).Value = 0.0745603636229859
$ws.Cells.Item(23, 5).Value = 2.301860511110718
$ws.Cells.Item(23, 6).Value = 0.7437029678328884
$ws.Cells.Item(23, 7).Value = 0.0007712498272074986
$ws.Cells.Item(23, 8).Value = 0.03518443906473134
$ws.Cells.Item(23, 9).Value = 0.01990882957891671
$ws.Cells.Item(23, 15).Value = 0
$ws.Cells.Item(23, 16).Value = 0
$ws.Cells.Item(23, 17).Value = 1.959238721746175

$ws.Cells.Item(24, 2).Value = 5.022345505370311
$ws.Cells.Item(24, 3).Value = 1.452387582760707
$ws.Cells.Item(24, 4).Value = 0.0656626206233426
$ws.Cells.Item(24, 5).Value = 1.967538615152051
$ws.Cells.Item(24, 6).Value = 0.6548386594035378
$ws.Cells.Item(24, 7).Value = 0.0007774558731151467
$ws.Cells.Item(24, 8).Value = 0.02548343545092102
$ws.Cells.Item(24, 9).Value = 0.01297716825878936
$ws.Cells.Item(24, 15).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 17).Value = 1.754635322995938

$ws.Cells.Item(25, 2).Value = 4.125615602166818
$ws.Cells.Item(25, 3).Value = 1.210872665251486
$ws.Cells.Item(25, 4).Value = 0.05608168448596018
$ws.Cells.Item(25, 5).Value = 1.611287029704712
$ws.Cells.Item(25, 6).Value = 0.5640920550256467
$ws.Cells.Item(25, 7).Value = 0.0007843738772956623
$ws.Cells.Item(25, 8).Value = 0.01649207729639268
$ws.Cells.Item(25, 9).Value = 0.007186637996994705
$ws.Cells.Item(25, 15).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 17).Value = 1.549647283128024
